$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unmerge existing merged header cells before restructuring
$ws.Range("A1:B1").UnMerge()
$ws.Range("C1:F1").UnMerge()

# Insert two new columns at C and D; this shifts old C,D,E,F -> E,F,G,H
$ws.Range("C:D").Insert()

# Move the "Moedas" header text from (now) C1 to E1
$ws.Range("C1").Value = $null
$ws.Range("E1").Value = "Moedas"

# New header row 2 values for inserted columns
$ws.Range("C2").Value = "Cotação da ação"
$ws.Range("D2").Value = "Valor atual total ação"

# New data values for inserted columns (rows 3-5) - store as text (leading
# apostrophe) to match the sibling columns (e.g. G/H) which hold
# numeric-looking strings, without introducing a new number-format style.
$ws.Range("C3").Value = "'1"
$ws.Range("D3").Value = "'3"

$ws.Range("C4").Value = "'1"
$ws.Range("D4").Value = "'3"

$ws.Range("C5").Value = "'1"
$ws.Range("D5").Value = "'3"

# Re-merge header cells over the new layout
$ws.Range("A1:D1").Merge()
$ws.Range("E1:H1").Merge()
